$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: new bibliography entry (MI-DRAGON) ---
$ws.Range("C20").Value = "3D MI-DRAGON:new model for the reconstruction of US FDA drug- target network and theoretical-experimental studies of inhibitors of rasagiline derivatives for AChE"
$ws.Range("D20").Value = 2012
$ws.Range("E20").Value = "desenvolvimento de modelo 3D mt-QSAR para a predição de molecuals inhibidoras usando MI-DRAGON (MARCH-INSIDE e DRAGON) para calcular os descritores 3D para os compostos da DrugBank database. Foram usados algoritmos ANN para determinar o melhor modelo no linear. o modelo obteve 87,03% de exatidão, sensibilidade de 85,35% e seletividade de 87,48%"
$ws.Range("F20").Value = "MLP(profile 37:37-24-1:1)"
$ws.Range("G20").Value = "DrugBank"
$ws.Range("H20").Value = "MI-DRAGON"
$ws.Range("I20").Value = "Prado-Prado, F.; García-Mera, X.; Escobar, M.; Alonso, N.;`nCaamaño, O.; Yañez, M.; González-Díaz, H. 3D MI-DRAGON:`nnew model for the reconstruction of US FDA drug- target network`nand theoretical-experimental studies of inhibitors of rasagiline derivatives`nfor AChE. Curr. Top. Med. Chem., 2012, 12(16), 1843-`n1865. [http://dx.doi.org/10.2174/1568026611209061843] [PMID:`n23030618]"

# G20 gets a distinct "right border only" look (picks up a new direct style
# in the saved workbook, same font/alignment as the rest of the row).
$g20 = $ws.Range("G20")
$g20.Borders.LineStyle = -4142
$g20.Borders.Item(10).LineStyle = 1
$g20.Borders.Item(10).Weight = 2
$g20.Interior.Pattern = -4142

$ws.Rows.Item(20).RowHeight = 92.4

# --- Row 21: fill in the missing "#Descritores" cell ---
$ws.Range("H21").Value = "3D-pharmacoporic"

# --- Update selection to match where the author ended up ---
$ws.Range("H21").Select()
